$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

$wsALC.Cells.Item(64, 8).Value = 3100.5625
$wsALC.Cells.Item(64, 9).Value = 3140.75
$wsALC.Cells.Item(64, 10).Value = 2980
$wsALC.Cells.Item(64, 11).Value = 3140.75
$wsALC.Cells.Item(64, 12).Value = 2980
$wsALC.Cells.Item(64, 13).Value = -2892.75
$wsALC.Cells.Item(64, 14).Value = -3476
$wsALC.Cells.Item(67, 8).Value = 3100.5625
$wsALC.Cells.Item(67, 9).Value = 3140.75
$wsALC.Cells.Item(67, 10).Value = 2980
$wsALC.Cells.Item(67, 11).Value = 3140.75
$wsALC.Cells.Item(67, 12).Value = 2980
$wsALC.Cells.Item(67, 13).Value = -2282.75
$wsALC.Cells.Item(67, 14).Value = -4696
$wsALC.Cells.Item(76, 8).Value = 3022.7778
$wsALC.Cells.Item(76, 9).Value = 2800
$wsALC.Cells.Item(76, 10).Value = 3802.5
$wsALC.Cells.Item(76, 11).Value = 2800
$wsALC.Cells.Item(76, 12).Value = 3802.5
$wsALC.Cells.Item(76, 13).Value = -2485
$wsALC.Cells.Item(76, 14).Value = -4432.5
$wsALC.Cells.Item(79, 8).Value = 3022.7778
$wsALC.Cells.Item(79, 9).Value = 2800
$wsALC.Cells.Item(79, 10).Value = 3802.5
$wsALC.Cells.Item(79, 11).Value = 2800
$wsALC.Cells.Item(79, 12).Value = 3802.5
$wsALC.Cells.Item(79, 13).Value = -1708
$wsALC.Cells.Item(79, 14).Value = -5986.5
$wsALC.Cells.Item(112, 8).Value = 10000961
$wsALC.Cells.Item(112, 10).Value = 10000961
$wsALC.Cells.Item(112, 12).Value = 30002883
$wsALC.Cells.Item(112, 14).Value = -30005099
$wsALC.Cells.Item(113, 8).Value = 3717.0715
$wsALC.Cells.Item(113, 9).Value = 3808.7778
$wsALC.Cells.Item(113, 10).Value = 3552
$wsALC.Cells.Item(113, 11).Value = 3808.7778
$wsALC.Cells.Item(113, 12).Value = 3552
$wsALC.Cells.Item(113, 13).Value = -554.7777999999998
$wsALC.Cells.Item(113, 14).Value = -10060
$wsALC.Cells.Item(137, 8).Value = 2131284.8
$wsALC.Cells.Item(137, 9).Value = 2634957.8
$wsALC.Cells.Item(137, 11).Value = 7904873.399999999
$wsALC.Cells.Item(137, 13).Value = -7902323.399999999
$wsALC.Cells.Item(138, 8).Value = 2426.5386
$wsALC.Cells.Item(138, 9).Value = 2084.7334
$wsALC.Cells.Item(138, 11).Value = 6254.2002
$wsALC.Cells.Item(138, 13).Value = -1114.2002
$wsALC.Cells.Item(141, 8).Value = 381353.3
$wsALC.Cells.Item(141, 10).Value = 564106.5600000001
$wsALC.Cells.Item(141, 12).Value = 1692319.68
$wsALC.Cells.Item(141, 14).Value = -1702679.68
$wsARM.Cells.Item(32, 8).Value = 6544.2046
$wsARM.Cells.Item(32, 9).Value = 5697.5303
$wsARM.Cells.Item(32, 10).Value = 20599
$wsARM.Cells.Item(32, 11).Value = 5697.5303
$wsARM.Cells.Item(32, 12).Value = 20599
$wsARM.Cells.Item(32, 13).Value = -5410.5303
$wsARM.Cells.Item(32, 14).Value = -21173
$wsARM.Cells.Item(110, 8).Value = 1246.7667
$wsARM.Cells.Item(110, 9).Value = 512.0833
$wsARM.Cells.Item(110, 11).Value = 512.0833
$wsARM.Cells.Item(110, 13).Value = 1532.9167
$wsBSM.Cells.Item(86, 8).Value = 23707.875
$wsBSM.Cells.Item(86, 9).Value = 2442.3333
$wsBSM.Cells.Item(86, 10).Value = 59150.445
$wsBSM.Cells.Item(86, 11).Value = 2442.3333
$wsBSM.Cells.Item(86, 12).Value = 59150.445
$wsBSM.Cells.Item(86, 13).Value = -1319.3333
$wsBSM.Cells.Item(86, 14).Value = -61396.445
$wsBSM.Cells.Item(89, 8).Value = 23707.875
$wsBSM.Cells.Item(89, 9).Value = 2442.3333
$wsBSM.Cells.Item(89, 10).Value = 59150.445
$wsBSM.Cells.Item(89, 11).Value = 12211.6665
$wsBSM.Cells.Item(89, 12).Value = 295752.225
$wsBSM.Cells.Item(89, 13).Value = -6595.666499999999
$wsBSM.Cells.Item(89, 14).Value = -306984.225
$wsBSM.Cells.Item(94, 8).Value = 724.36365
$wsBSM.Cells.Item(94, 9).Value = 762
$wsBSM.Cells.Item(94, 11).Value = 762
$wsBSM.Cells.Item(94, 13).Value = -311
$wsBSM.Cells.Item(99, 8).Value = 2159
$wsBSM.Cells.Item(99, 9).Value = 1830.5385
$wsBSM.Cells.Item(99, 10).Value = 2514.8333
$wsBSM.Cells.Item(99, 11).Value = 1830.5385
$wsBSM.Cells.Item(99, 12).Value = 2514.8333
$wsBSM.Cells.Item(99, 13).Value = -332.5385000000001
$wsBSM.Cells.Item(99, 14).Value = -5510.8333
$wsBSM.Cells.Item(105, 8).Value = 1472.5526
$wsBSM.Cells.Item(105, 9).Value = 1455.0952
$wsBSM.Cells.Item(105, 11).Value = 1455.0952
$wsBSM.Cells.Item(105, 13).Value = 291.9048
$wsCRP.Cells.Item(52, 8).Value = 34583.332
$wsCRP.Cells.Item(52, 10).Value = 34583.332
$wsCRP.Cells.Item(52, 12).Value = 34583.332
$wsCRP.Cells.Item(52, 14).Value = -35171.332
$wsCRP.Cells.Item(62, 8).Value = 4498.737
$wsCRP.Cells.Item(62, 9).Value = 3319.2856
$wsCRP.Cells.Item(62, 10).Value = 7801.2
$wsCRP.Cells.Item(62, 11).Value = 3319.2856
$wsCRP.Cells.Item(62, 12).Value = 7801.2
$wsCRP.Cells.Item(62, 13).Value = -2695.2856
$wsCRP.Cells.Item(62, 14).Value = -9049.200000000001
$wsCRP.Cells.Item(65, 8).Value = 4498.737
$wsCRP.Cells.Item(65, 9).Value = 3319.2856
$wsCRP.Cells.Item(65, 10).Value = 7801.2
$wsCRP.Cells.Item(65, 11).Value = 16596.428
$wsCRP.Cells.Item(65, 12).Value = 39006
$wsCRP.Cells.Item(65, 13).Value = -13476.428
$wsCRP.Cells.Item(65, 14).Value = -45246
$wsCRP.Cells.Item(99, 8).Value = 3086.875
$wsCRP.Cells.Item(99, 9).Value = 1099.1111
$wsCRP.Cells.Item(99, 11).Value = 1099.1111
$wsCRP.Cells.Item(99, 13).Value = 398.8888999999999
$wsCRP.Cells.Item(126, 8).Value = 3086.875
$wsCRP.Cells.Item(126, 9).Value = 1099.1111
$wsCRP.Cells.Item(126, 11).Value = 3297.3333
$wsCRP.Cells.Item(126, 13).Value = -827.3333000000002
$wsCUL.Cells.Item(131, 8).Value = 1177
$wsCUL.Cells.Item(131, 9).Value = 2185
$wsCUL.Cells.Item(131, 10).Value = 1069
$wsCUL.Cells.Item(131, 11).Value = 6555
$wsCUL.Cells.Item(131, 12).Value = 3207
$wsCUL.Cells.Item(131, 13).Value = -1515
$wsCUL.Cells.Item(131, 14).Value = -13287
$wsGSM.Cells.Item(19, 8).Value = 24301.2
$wsGSM.Cells.Item(19, 9).Value = 2750
$wsGSM.Cells.Item(19, 10).Value = 38668.668
$wsGSM.Cells.Item(19, 11).Value = 2750
$wsGSM.Cells.Item(19, 12).Value = 38668.668
$wsGSM.Cells.Item(19, 13).Value = -2462
$wsGSM.Cells.Item(19, 14).Value = -39244.668
$wsGSM.Cells.Item(80, 8).Value = 4056.5
$wsGSM.Cells.Item(80, 9).Value = 3058.6
$wsGSM.Cells.Item(80, 11).Value = 3058.6
$wsGSM.Cells.Item(80, 13).Value = -2060.6
$wsGSM.Cells.Item(83, 8).Value = 4056.5
$wsGSM.Cells.Item(83, 9).Value = 3058.6
$wsGSM.Cells.Item(83, 11).Value = 15293
$wsGSM.Cells.Item(83, 13).Value = -10301
$wsLTW.Cells.Item(32, 8).Value = 38495.168
$wsLTW.Cells.Item(32, 9).Value = 3642
$wsLTW.Cells.Item(32, 11).Value = 3642
$wsLTW.Cells.Item(32, 13).Value = -3325
$wsLTW.Cells.Item(46, 8).Value = 1292.1154
$wsLTW.Cells.Item(46, 9).Value = 960.7143
$wsLTW.Cells.Item(46, 10).Value = 2684
$wsLTW.Cells.Item(46, 11).Value = 960.7143
$wsLTW.Cells.Item(46, 12).Value = 2684
$wsLTW.Cells.Item(46, 13).Value = -772.7143
$wsLTW.Cells.Item(46, 14).Value = -3060
$wsLTW.Cells.Item(68, 8).Value = 1791.3043
$wsLTW.Cells.Item(68, 9).Value = 1020
$wsLTW.Cells.Item(68, 10).Value = 6933.3335
$wsLTW.Cells.Item(68, 11).Value = 1020
$wsLTW.Cells.Item(68, 12).Value = 6933.3335
$wsLTW.Cells.Item(68, 13).Value = -271
$wsLTW.Cells.Item(68, 14).Value = -8431.333500000001
$wsLTW.Cells.Item(71, 8).Value = 1791.3043
$wsLTW.Cells.Item(71, 9).Value = 1020
$wsLTW.Cells.Item(71, 10).Value = 6933.3335
$wsLTW.Cells.Item(71, 11).Value = 5100
$wsLTW.Cells.Item(71, 12).Value = 34666.6675
$wsLTW.Cells.Item(71, 13).Value = -1356
$wsLTW.Cells.Item(71, 14).Value = -42154.6675
$wsLTW.Cells.Item(82, 8).Value = 3100.1538
$wsLTW.Cells.Item(82, 9).Value = 2050.25
$wsLTW.Cells.Item(82, 11).Value = 2050.25
$wsLTW.Cells.Item(82, 13).Value = -1689.25
$wsLTW.Cells.Item(85, 8).Value = 3100.1538
$wsLTW.Cells.Item(85, 9).Value = 2050.25
$wsLTW.Cells.Item(85, 11).Value = 2050.25
$wsLTW.Cells.Item(85, 13).Value = -802.25
$wsWVR.Cells.Item(62, 8).Value = 4835.846
$wsWVR.Cells.Item(62, 9).Value = 4475.5
$wsWVR.Cells.Item(62, 11).Value = 4475.5
$wsWVR.Cells.Item(62, 13).Value = -3851.5
$wsWVR.Cells.Item(65, 8).Value = 4835.846
$wsWVR.Cells.Item(65, 9).Value = 4475.5
$wsWVR.Cells.Item(65, 11).Value = 22377.5
$wsWVR.Cells.Item(65, 13).Value = -19257.5
$wsWVR.Cells.Item(122, 8).Value = 287668.84
$wsWVR.Cells.Item(122, 9).Value = 527864.2
$wsWVR.Cells.Item(122, 10).Value = 2436.8125
$wsWVR.Cells.Item(122, 11).Value = 1583592.6
$wsWVR.Cells.Item(122, 12).Value = 7310.4375
$wsWVR.Cells.Item(122, 13).Value = -1581142.6
$wsWVR.Cells.Item(122, 14).Value = -12210.4375
$wsWVR.Cells.Item(126, 8).Value = 3032330
$wsWVR.Cells.Item(126, 9).Value = 1381.0454
$wsWVR.Cells.Item(126, 11).Value = 4143.1362
$wsWVR.Cells.Item(126, 13).Value = -1673.1362
